$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '58.100.62'
$c.ClearFormats()
$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  -1.90%  '
$c.ClearFormats()
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.289.44'
$c.ClearFormats()
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  -2.18%  '
$c.ClearFormats()
$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  -0.17%  '
$c.ClearFormats()
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '531.27'
$c.ClearFormats()
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  -4.89%  '
$c.ClearFormats()
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '131.03'
$c.ClearFormats()
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  -0.60%  '
$c.ClearFormats()
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '
$c.ClearFormats()
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.582'
$c.ClearFormats()
$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +1.08%  '
$c.ClearFormats()
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '2.288.70'
$c.ClearFormats()
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  -1.95%  '
$c.ClearFormats()
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0995'
$c.ClearFormats()
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  -3.94%  '
$c.ClearFormats()
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '5.47'
$c.ClearFormats()
$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  -2.05%  '
$c.ClearFormats()
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '
$c.ClearFormats()
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.329'
$c.ClearFormats()
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -2.78%  '
$c.ClearFormats()
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '23.38'
$c.ClearFormats()
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  -2.27%  '
$c.ClearFormats()
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '2.695.46'
$c.ClearFormats()
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -2.46%  '
$c.ClearFormats()
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '57.991.12'
$c.ClearFormats()
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  -2.13%  '
$c.ClearFormats()
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.0000132'
$c.ClearFormats()
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  -2.65%  '
$c.ClearFormats()
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '2.268.37'
$c.ClearFormats()
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  -3.20%  '
$c.ClearFormats()
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '10.50'
$c.ClearFormats()
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  -3.79%  '
$c.ClearFormats()
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '4.17'
$c.ClearFormats()
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  -5.56%  '
$c.ClearFormats()
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '311.58'
$c.ClearFormats()
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  -2.21%  '
$c.ClearFormats()
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '6.40'
$c.ClearFormats()
$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -2.54%  '
$c.ClearFormats()
$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +0.06%  '
$c.ClearFormats()
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '62.33'
$c.ClearFormats()
$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  -2.13%  '
$c.ClearFormats()
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '0.167'
$c.ClearFormats()
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  -1.60%  '
$c.ClearFormats()
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  -0.20%  '
$c.ClearFormats()
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '8.01'
$c.ClearFormats()
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  -3.88%  '
$c.ClearFormats()
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '1.26'
$c.ClearFormats()
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -6.30%  '
$c.ClearFormats()
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '170.08'
$c.ClearFormats()
$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  -0.50%  '
$c.ClearFormats()
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '1.70'
$c.ClearFormats()
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  -5.43%  '
$c.ClearFormats()
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.0₃0716'
$c.ClearFormats()
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -3.41%  '
$c.ClearFormats()
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '5.72'
$c.ClearFormats()
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -3.25%  '
$c.ClearFormats()
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  -5.09%  '
$c.ClearFormats()
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.377'
$c.ClearFormats()
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  -4.62%  '
$c.ClearFormats()
$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  +0.00%  '
$c.ClearFormats()
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '17.74'
$c.ClearFormats()
$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -1.14%  '
$c.ClearFormats()
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -0.21%  '
$c.ClearFormats()
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '1.23'
$c.ClearFormats()
$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -5.21%  '
$c.ClearFormats()
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '3.89'
$c.ClearFormats()
$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -4.12%  '
$c.ClearFormats()
$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  -0.16%  '
$c.ClearFormats()
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '1.49'
$c.ClearFormats()
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  -4.50%  '
$c.ClearFormats()
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '140.47'
$c.ClearFormats()
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  -2.63%  '
$c.ClearFormats()
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '287.08'
$c.ClearFormats()
$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -7.43%  '
$c.ClearFormats()
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '3.41'
$c.ClearFormats()
$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  -1.77%  '
$c.ClearFormats()
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.0947'
$c.ClearFormats()
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -0.67%  '
$c.ClearFormats()
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.0494'
$c.ClearFormats()
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -1.58%  '
$c.ClearFormats()
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.552'
$c.ClearFormats()
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -1.84%  '
$c.ClearFormats()
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '18.03'
$c.ClearFormats()
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -3.46%  '
$c.ClearFormats()
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.0209'
$c.ClearFormats()
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  -2.71%  '
$c.ClearFormats()
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -1.20%  '
$c.ClearFormats()
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -0.74%  '
$c.ClearFormats()
